$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift rows 2-19 up by two (copy row n+2 contents for columns C-H into row n)
# Done in increasing row order so we always read original (not-yet-overwritten) data,
# since the source row (n+2) is only overwritten by this same loop once n reaches n+2,
# which happens strictly later (n+2 > n).
for ($n = 2; $n -le 19; $n++) {
    $src = $n + 2
    for ($col = 3; $col -le 8; $col++) {
        $val = $ws.Cells.Item($src, $col).Value()
        $ws.Cells.Item($n, $col).Value = $val
    }
}

# Step 2: new data for rows 20-31 (A: timestamp, B: label "struggle", C-H: sensor data)
$newRows = @{
    20 = @(1800, -3.367526054382324, 0.915987193584442, -1.626443386077881, -1.0144944190979, -1.167210817337036, 0.6551529765129089)
    21 = @(1900, -1.118759155273438, 2.792432069778442, -4.963344097137451, -0.3874412775039673, -0.7050912380218506, 0.0914770737290382)
    22 = @(2000, -1.557756900787354, -0.5582034587860107, -0.2619988918304443, -0.1458440721035003, -0.3762930035591125, -0.0704022198915481)
    23 = @(2100, 1.715949058532715, -1.576748490333557, 5.096891403198242, 0.2157881408929824, 0.3240640163421631, 0.0951422601938247)
    24 = @(2200, -0.8243503570556641, 0.5943599939346313, 1.927432060241699, 0.0403171069920063, 0.1484402567148208, -0.0852157026529312)
    25 = @(2300, -0.2981023788452148, 1.024843096733093, 0.8517363667488098, 0.1014036312699318, 0.3179553747177124, 0.0390953756868839)
    26 = @(2400, 0.4371089935302734, 0.3337190449237823, -0.154114544391632, 0.052381694316864, 0.1099557429552078, 0.0681114718317985)
    27 = @(2500, -0.5641984939575195, -0.3292053341865539, -0.326197862625122, 0.0522289797663688, -0.4196644127368927, 0.2273945808410644)
    28 = @(2600, 0.131052017211914, 0.5107872486114502, 0.0540084838867187, 0.0937678143382072, -0.1565342247486114, 0.0675006061792373)
    29 = @(2700, -0.1777238845825195, 0.4102384448051452, 0.1352127194404602, -0.0591012127697467, 0.0331394411623477, 0.0291688162833452)
    30 = @(2800, 0.2971744537353515, 0.0503720641136169, -0.09675005078315729, -0.0047342055477201, 0.0655152946710586, -0.0163406450301408)
    31 = @(2900, -0.1711950302124023, -0.0577961653470993, 0.0454300940036773, -0.00534507073462, 0.0612392425537109, -0.0022907445672899)
}

foreach ($r in $newRows.Keys) {
    $row = $newRows[$r]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "struggle"
    for ($i = 1; $i -le 6; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $row[$i]
    }
}

Write-Host "Edit complete. New dimension: " $ws.UsedRange.Address()